$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove cells that no longer have forecast values (bug fix removed erroneous early entries)
$ws.Range("C2").Value = $null
$ws.Range("E2").Value = $null
$ws.Range("C3").Value = $null
$ws.Range("C4").Value = $null

# Update values with corrected floating point forecaster output
$ws.Range("E4").Value = 1.087227286828263
$ws.Range("C5").Value = -5.478010998490146
$ws.Range("C6").Value = -4.774178217057756
$ws.Range("E6").Value = -0.8523446516643496
$ws.Range("E7").Value = -1.305195642355672
$ws.Range("E8").Value = 1.390521443873438
$ws.Range("C9").Value = 3.371423250978833
$ws.Range("E9").Value = 0.806063216063202
$ws.Range("E11").Value = 2.76635821344573
$ws.Range("C12").Value = 1.239479831392831
$ws.Range("C13").Value = -0.03183655677960751
$ws.Range("E13").Value = 1.102200073559856
$ws.Range("C14").Value = 0.2379616621360992
$ws.Range("C15").Value = 1.812248956008777
$ws.Range("E15").Value = 1.209672013646301
$ws.Range("E16").Value = 0.6176326357195894
$ws.Range("C18").Value = 1.470039379455734
$ws.Range("E18").Value = 1.577608035818323
$ws.Range("C19").Value = 1.575690123464613
$ws.Range("E19").Value = 1.643656926428538
$ws.Range("C20").Value = 1.638797242243228
$ws.Range("E20").Value = 1.369334405341593
$ws.Range("C23").Value = 2.507284186438108
$ws.Range("E23").Value = 2.112314908467128
$ws.Range("C24").Value = 2.337818484846443
$ws.Range("C25").Value = 2.354760705778181
$ws.Range("E25").Value = 2.107524645430892
$ws.Range("E26").Value = 1.36203066512679
$ws.Range("C27").Value = 0.803755999809086
$ws.Range("E27").Value = 1.649904670037827
$ws.Range("C28").Value = 0.8311911554373275
$ws.Range("E28").Value = 1.758956425699276
$ws.Range("C29").Value = 0.7024402883234249
$ws.Range("E31").Value = 0.9721240557711175
$ws.Range("C32").Value = -1.538034740964356
$ws.Range("E32").Value = -0.7351085756681197
$ws.Range("C33").Value = -6.356537224117531
$ws.Range("E33").Value = -8.821046965146573
$ws.Range("E34").Value = -0.1895486537906388
$ws.Range("C35").Value = 0.4146413258694359
$ws.Range("E35").Value = -0.8956136585515861
$ws.Range("C36").Value = -0.5490727792360039
$ws.Range("E37").Value = -0.454653018564577
$ws.Range("C38").Value = 1.099928004397577
$ws.Range("C40").Value = 1.687572871803722
$ws.Range("C41").Value = 2.234093617591992
$ws.Range("C42").Value = 2.310042359896247
$ws.Range("E43").Value = 4.881781055849221
$ws.Range("C44").Value = -0.02880469535951891
$ws.Range("E46").Value = 1.063035646777677
$ws.Range("E47").Value = 0.8017085309184768
$ws.Range("C48").Value = -0.1645795020818963
$ws.Range("E48").Value = 0.2574142441027716
$ws.Range("C49").Value = -0.2058547204034422
$ws.Range("E49").Value = 0.2971745009357374
$ws.Range("C50").Value = -0.3101476031197037
$ws.Range("C51").Value = -0.2675654179851272
$ws.Range("E51").Value = -0.1518868483142199
$ws.Range("C52").Value = 0.06409464788890151
$ws.Range("E52").Value = -0.06354501920062816
$ws.Range("C53").Value = -0.1271672627326415
$ws.Range("E53").Value = -0.1556121492501283
